$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = "[4.075321314796454, 9.505886811516028]"
$ws.Range("M2").Value = [double]"1.366120800039283e-06"
$ws.Range("N2").Value = [double]"1.366120800039283e-06"
$ws.Range("P2").Value = "[-1.54721079636254, -0.7170001251436169]"
$ws.Range("Q2").Value = [double]"1.519782089065558e-07"
$ws.Range("R2").Value = [double]"1.519782089065558e-07"
$ws.Range("T2").Value = "[7.706942433375938, 10.591590809488903]"
$ws.Range("X2").Value = [double]"2.87909909909916"
$ws.Range("Y2").Value = [double]"6.212792792792916"

# Row 3 updates
$ws.Range("L3").Value = "[4.4433339711680695, 9.204601443054532]"
$ws.Range("M3").Value = [double]"3.191533193280804e-08"
$ws.Range("N3").Value = [double]"6.383066386561609e-08"
$ws.Range("P3").Value = "[2.207605648468503, 2.937184723176043]"
$ws.Range("T3").Value = "[7.807052220851402, 10.32755714032704]"
$ws.Range("X3").Value = [double]"12.73817817817834"
$ws.Range("Y3").Value = [double]"15.51567567567587"

Write-Host "Edits applied"
